$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "1= Male, 0= Female"

$ws.Range("A10").Value = "Waist circumference (in)"
$ws.Range("B10").Value = "Waist_Circumference"

$ws.Range("A11").Value = "BMI (kg/m^2)"
$ws.Range("B11").Value = "BMI"

$ws.Range("B12").Value = "BMI_Range"
$ws.Range("C12").Value = "Text"
$ws.Range("D12").Value = "Underweight, Healthy Weight, Overweight, Obese"
$ws.Range("A12").Value = "BMI Range, based on BMI growth curves normed for children, accounting for age and sex"

$ws.Range("A13").Select()
